$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal text value, avoiding Excel's automatic
# type inference (date / number parsing) that would otherwise change the
# stored cell type away from a plain shared string.
function Set-TextValue($row, $col, [string]$text) {
    $cell = $ws.Cells.Item($row, $col)
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
    $excel.CutCopyMode = $false
}

$rowData = @(
    @{row=2; A=1; C="Injuries"; D=35.07438; E=-85.17346; F="2018-12-10"; G="22:50:19"; H="6317-6339 BONNY OAKS DR"; K="CHATTANOOGA"; P="22"}
    @{row=3; A=18; C="Injuries"; D=35.068625; E=-85.263132; F="2018-12-10"; G="18:23:43"; H="Riverside Dr / Latta St"; K="CHATTANOOGA"; P="18"}
    @{row=4; A=24; C="Injuries"; D=35.027174; E=-85.190759; F="2018-12-10"; G="17:54:52"; H="6100 Lee Hwy"; K="CHATTANOOGA"; P="17"}
    @{row=5; A=25; C="Injuries"; D=35.027174; E=-85.190759; F="2018-12-10"; G="17:54:15"; H="6100 Lee Hwy"; K="CHATTANOOGA"; P="17"}
    @{row=6; A=32; C="Injuries"; D=35.087611; E=-85.214113; F="2018-12-10"; G="17:14:49"; H="520 - 549 Highway 153 Nb"; K="CHATTANOOGA"; P="17"}
    @{row=7; A=33; C="Injuries"; D=35.087611; E=-85.214113; F="2018-12-10"; G="17:14:18"; H="520 HIGHWAY 153 NB"; K="CHATTANOOGA"; P="17"}
    @{row=8; A=37; C="Unknown Injuries"; D=34.996398; E=-85.245409; F="2018-12-10"; G="16:58:39"; H="4150 Ringgold Rd"; K="EAST RIDGE"; P="16"}
    @{row=9; A=46; C="Injuries"; D=35.201065; E=-85.238303; F="2018-12-10"; G="16:06:43"; H="7600 DAYTON PIKE"; K="SODDY DAISY"; P="16"}
    @{row=10; A=46; C="Injuries"; D=35.201065; E=-85.238303; F="2018-12-10"; G="16:06:43"; H="7600 DAYTON PIKE"; K="SODDY DAISY"; P="16"}
    @{row=11; A=48; C="Entrapment"; D=35.201065; E=-85.238303; F="2018-12-10"; G="16:05:33"; H="7600 DAYTON PIKE"; K="SODDY DAISY"; P="16"}
    @{row=12; A=49; C="Entrapment"; D=35.201065; E=-85.238303; F="2018-12-10"; G="16:04:58"; H="7600 DAYTON PIKE"; K="SODDY DAISY"; P="16"}
    @{row=13; A=51; C="Injuries"; D=35.12002; E=-85.249142; F="2018-12-10"; G="16:04:49"; H="4513 HIXSON PIKE"; K="CHATTANOOGA"; P="16"}
    @{row=14; A=77; C="Injuries"; D=35.24122; E=-85.175514; F="2018-12-10"; G="11:23:51"; H="600 GREEN POND RD"; K="SODDY DAISY"; P="11"}
    @{row=15; A=81; C="Injuries"; D=35.084112; E=-85.208627; F="2018-12-10"; G="10:27:18"; H="480 HIGHWAY 153 SB"; K="CHATTANOOGA"; P="10"}
    @{row=16; A=82; C="Injuries"; D=35.084112; E=-85.208627; F="2018-12-10"; G="10:27:18"; H="480 HIGHWAY 153 SB"; K="CHATTANOOGA"; P="10"}
)

foreach ($d in $rowData) {
    $r = $d.row
    $ws.Cells.Item($r, 1).Value = $d.A
    $ws.Cells.Item($r, 3).Value = $d.C
    $ws.Cells.Item($r, 4).Value = $d.D
    $ws.Cells.Item($r, 5).Value = $d.E
    Set-TextValue $r 6 $d.F
    $ws.Cells.Item($r, 7).Value = $d.G
    $ws.Cells.Item($r, 8).Value = $d.H
    $ws.Cells.Item($r, 11).Value = $d.K
    Set-TextValue $r 16 $d.P
    $ws.Cells.Item($r, 25).Value = 12
}

# New rows 11-16 need column A formatted like the existing data rows (s="1").
# Copy the number format from an existing styled cell (A2) onto the new ones.
$ws.Cells.Item(2, 1).Copy()
$ws.Range($ws.Cells.Item(11, 1), $ws.Cells.Item(16, 1)).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

